# Se modificó el tiempo de respuesta del sensor 2.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Sensor 2 (TMP36) row is row 5; column D holds "Tiempo de respuesta (100%)".
$ws.Range("D5").Value = "8 min"

# Leftover formatting on an adjacent, otherwise-empty cell (underline font,
# no border/fill) - matches the extra style picked up while editing row 5.
$ws.Range("P5").Font.Underline = $true

# Leave the selection/scroll position where the edit happened.
$excel.ActiveWindow.ScrollColumn = 9
$excel.ActiveWindow.ScrollRow = 3
$ws.Range("P5").Select()
